$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H17").Value = 674.19354
$ws.Range("J17").Value = 674.19354
$ws.Range("L17").Value = 2022.58062
$ws.Range("N17").Value = -2358.58062

$ws.Range("H57").Value = 65533.332
$ws.Range("J57").Value = 65533.332
$ws.Range("L57").Value = 196599.996
$ws.Range("N57").Value = -197597.996

$ws.Range("H98").Value = 6225.448
$ws.Range("I98").Value = 6542.148
$ws.Range("J98").Value = 1950
$ws.Range("K98").Value = 6542.148
$ws.Range("L98").Value = 1950
$ws.Range("M98").Value = -5044.148
$ws.Range("N98").Value = -4946

$ws.Range("H122").Value = 6225.448
$ws.Range("I122").Value = 6542.148
$ws.Range("J122").Value = 1950
$ws.Range("K122").Value = 19626.444
$ws.Range("L122").Value = 5850
$ws.Range("M122").Value = -17176.444
$ws.Range("N122").Value = -10750

$ws.Range("H132").Value = 4390135.5
$ws.Range("I132").Value = 5467011.5
$ws.Range("K132").Value = 16401034.5
$ws.Range("M132").Value = -16398504.5

$ws.Range("H138").Value = 1501.81
$ws.Range("I138").Value = 814.7857
$ws.Range("J138").Value = 1613.6511
$ws.Range("K138").Value = 2444.3571
$ws.Range("L138").Value = 4840.9533
$ws.Range("M138").Value = 2695.6429
$ws.Range("N138").Value = -15120.9533

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 2688.3684
$ws.Range("I32").Value = 2456.973
$ws.Range("K32").Value = 2456.973
$ws.Range("M32").Value = -2169.973

$ws.Range("H45").Value = 1309.2667
$ws.Range("I45").Value = 1276
$ws.Range("J45").Value = 1442.3334
$ws.Range("K45").Value = 1276
$ws.Range("L45").Value = 1442.3334
$ws.Range("M45").Value = -899
$ws.Range("N45").Value = -2196.3334

$ws.Range("H61").Value = 1133.8628
$ws.Range("I61").Value = 1037.7693
$ws.Range("J61").Value = 1446.1666
$ws.Range("K61").Value = 1037.7693
$ws.Range("L61").Value = 1446.1666
$ws.Range("M61").Value = -825.7692999999999
$ws.Range("N61").Value = -1870.1666

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 902.9167
$ws.Range("I74").Value = 497.37036
$ws.Range("J74").Value = 2119.5557
$ws.Range("K74").Value = 497.37036
$ws.Range("L74").Value = 2119.5557
$ws.Range("M74").Value = 376.62964
$ws.Range("N74").Value = -3867.5557

$ws.Range("H77").Value = 902.9167
$ws.Range("I77").Value = 497.37036
$ws.Range("J77").Value = 2119.5557
$ws.Range("K77").Value = 2486.8518
$ws.Range("L77").Value = 10597.7785
$ws.Range("M77").Value = 1881.1482
$ws.Range("N77").Value = -19333.7785

$ws.Range("H110").Value = 1638.1904
$ws.Range("I110").Value = 1239.3334
$ws.Range("J110").Value = 2635.3333
$ws.Range("K110").Value = 1239.3334
$ws.Range("L110").Value = 2635.3333
$ws.Range("M110").Value = 805.6666
$ws.Range("N110").Value = -6725.3333

$ws.Range("H122").Value = 1295.76
$ws.Range("I122").Value = 1140
$ws.Range("J122").Value = 2438
$ws.Range("K122").Value = 3420
$ws.Range("L122").Value = 7314
$ws.Range("M122").Value = -970
$ws.Range("N122").Value = -12214

$ws.Range("H132").Value = 1324.8448
$ws.Range("I132").Value = 1044.7317
$ws.Range("K132").Value = 3134.1951
$ws.Range("M132").Value = -604.1950999999999

$ws.Range("H136").Value = 1133.8628
$ws.Range("I136").Value = 1037.7693
$ws.Range("J136").Value = 1446.1666
$ws.Range("K136").Value = 3113.3079
$ws.Range("L136").Value = 4338.4998
$ws.Range("M136").Value = -563.3078999999998
$ws.Range("N136").Value = -9438.4998

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H99").Value = 71429870
$ws.Range("I99").Value = 83334456
$ws.Range("K99").Value = 83334456
$ws.Range("M99").Value = -83332958

$ws.Range("H107").Value = 1785.1578
$ws.Range("I107").Value = 1471.8572
$ws.Range("J107").Value = 2662.4
$ws.Range("K107").Value = 1471.8572
$ws.Range("L107").Value = 2662.4
$ws.Range("M107").Value = 448.1428000000001
$ws.Range("N107").Value = -6502.4

$ws.Range("H134").Value = 3932.558
$ws.Range("I134").Value = 999.7879
$ws.Range("J134").Value = 13610.7
$ws.Range("K134").Value = 2999.3637
$ws.Range("L134").Value = 40832.10000000001
$ws.Range("M134").Value = -464.3636999999999
$ws.Range("N134").Value = -45902.10000000001

$ws.Range("H141").Value = 74998
$ws.Range("J141").Value = 74998
$ws.Range("L141").Value = 74998
$ws.Range("N141").Value = -85358

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 1702.2285
$ws.Range("I31").Value = 1791.2963
$ws.Range("J31").Value = 1401.625
$ws.Range("K31").Value = 1791.2963
$ws.Range("L31").Value = 1401.625
$ws.Range("M31").Value = -1496.2963
$ws.Range("N31").Value = -1991.625

$ws.Range("H34").Value = 1702.2285
$ws.Range("I34").Value = 1791.2963
$ws.Range("J34").Value = 1401.625
$ws.Range("K34").Value = 1791.2963
$ws.Range("L34").Value = 1401.625
$ws.Range("M34").Value = -1589.2963
$ws.Range("N34").Value = -1805.625

$ws.Range("H51").Value = 17333.334
$ws.Range("J51").Value = 19800
$ws.Range("L51").Value = 19800
$ws.Range("N51").Value = -21272

$ws.Range("H56").Value = 20103
$ws.Range("J56").Value = 20103
$ws.Range("L56").Value = 20103
$ws.Range("N56").Value = -21793

$ws.Range("H58").Value = 694.2727
$ws.Range("I58").Value = 616.5
$ws.Range("J58").Value = 958.7
$ws.Range("K58").Value = 616.5
$ws.Range("L58").Value = 958.7
$ws.Range("M58").Value = -413.5
$ws.Range("N58").Value = -1364.7

$ws.Range("H59").Value = 26000
$ws.Range("J59").Value = 26000
$ws.Range("L59").Value = 26000
$ws.Range("N59").Value = -28290

$ws.Range("H61").Value = 17333.334
$ws.Range("J61").Value = 19800
$ws.Range("L61").Value = 19800
$ws.Range("N61").Value = -20496

$ws.Range("H132").Value = 3408.7036
$ws.Range("I132").Value = 3635.6743
$ws.Range("K132").Value = 10907.0229
$ws.Range("M132").Value = -8377.0229

$ws.Range("H134").Value = 889.9091
$ws.Range("I134").Value = 918.08887
$ws.Range("J134").Value = 763.1
$ws.Range("K134").Value = 2754.26661
$ws.Range("L134").Value = 2289.3
$ws.Range("M134").Value = -219.2666100000001
$ws.Range("N134").Value = -7359.3

$ws.Range("H136").Value = 694.2727
$ws.Range("I136").Value = 616.5
$ws.Range("J136").Value = 958.7
$ws.Range("K136").Value = 1849.5
$ws.Range("L136").Value = 2876.1
$ws.Range("M136").Value = 700.5
$ws.Range("N136").Value = -7976.1

$ws.Range("H139").Value = 44780
$ws.Range("J139").Value = 44780
$ws.Range("L139").Value = 44780
$ws.Range("N139").Value = -55060

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H5").Value = 1384.4
$ws.Range("I5").Value = 1515.5555
$ws.Range("J5").Value = 941.75
$ws.Range("K5").Value = 4546.666499999999
$ws.Range("L5").Value = 2825.25
$ws.Range("M5").Value = -4434.666499999999
$ws.Range("N5").Value = -3049.25

$ws.Range("H14").Value = 699.5
$ws.Range("I14").Value = 699.5
$ws.Range("K14").Value = 2098.5
$ws.Range("M14").Value = -1925.5

$ws.Range("H135").Value = 1384.4
$ws.Range("I135").Value = 1515.5555
$ws.Range("J135").Value = 941.75
$ws.Range("K135").Value = 13639.9995
$ws.Range("L135").Value = 8475.75
$ws.Range("M135").Value = -11104.9995
$ws.Range("N135").Value = -13545.75

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H102").Value = 6941.5557
$ws.Range("I102").Value = 5611.385
$ws.Range("J102").Value = 10400
$ws.Range("K102").Value = 5611.385
$ws.Range("L102").Value = 10400
$ws.Range("M102").Value = -3989.385
$ws.Range("N102").Value = -13644

$ws.Range("H122").Value = 1831.3182
$ws.Range("I122").Value = 1157.625
$ws.Range("J122").Value = 3627.8333
$ws.Range("K122").Value = 3472.875
$ws.Range("L122").Value = 10883.4999
$ws.Range("M122").Value = -1022.875
$ws.Range("N122").Value = -15783.4999

$ws.Range("H126").Value = 2076.75
$ws.Range("I126").Value = 1760
$ws.Range("K126").Value = 5280
$ws.Range("M126").Value = -2810

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H22").Value = 1020.2
$ws.Range("I22").Value = 733
$ws.Range("K22").Value = 733
$ws.Range("M22").Value = -438

$ws.Range("H27").Value = 1020.2
$ws.Range("I27").Value = 733
$ws.Range("K27").Value = 733
$ws.Range("M27").Value = -626

$ws.Range("H122").Value = 35731270
$ws.Range("I122").Value = 35731270
$ws.Range("K122").Value = 107193810
$ws.Range("M122").Value = -107191360

$ws.Range("H136").Value = 7965.9375
$ws.Range("I136").Value = 12511.111
$ws.Range("J136").Value = 2122.1428
$ws.Range("K136").Value = 37533.333
$ws.Range("L136").Value = 6366.428400000001
$ws.Range("M136").Value = -34983.333
$ws.Range("N136").Value = -11466.4284

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H32").Value = 9900
$ws.Range("I32").Value = 9900
$ws.Range("K32").Value = 9900
$ws.Range("M32").Value = -9583

$ws.Range("H113").Value = 407.8125
$ws.Range("I113").Value = 286.22223
$ws.Range("J113").Value = 564.1429000000001
$ws.Range("K113").Value = 858.66669
$ws.Range("L113").Value = 1692.4287
$ws.Range("M113").Value = 1311.33331
$ws.Range("N113").Value = -6032.4287

$ws.Range("H127").Value = 75000
$ws.Range("J127").Value = 75000
$ws.Range("L127").Value = 75000
$ws.Range("N127").Value = -84920

$ws.Range("H132").Value = 1869.8363
$ws.Range("I132").Value = 1944.6666
$ws.Range("K132").Value = 5833.9998
$ws.Range("M132").Value = -3303.9998

$ws.Range("H136").Value = 478.84848
$ws.Range("I136").Value = 326.54166
$ws.Range("J136").Value = 885
$ws.Range("K136").Value = 979.6249799999999
$ws.Range("L136").Value = 2655
$ws.Range("M136").Value = 1570.37502
$ws.Range("N136").Value = -7755
